$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 18.76993966666667
$ws.Cells.Item(2, 8).Value = 56.309819
$ws.Cells.Item(2, 9).Value = 0.1007685501185251
$ws.Cells.Item(2, 10).Value = 0.1007685501185251
$ws.Cells.Item(2, 13).Value = 24.41295733333333
$ws.Cells.Item(2, 14).Value = 73.238872
$ws.Cells.Item(2, 15).Value = 0.1430960671192788
$ws.Cells.Item(2, 16).Value = 0.1430960671192788
$ws.Cells.Item(2, 17).Value = 458.2297362315742
$ws.Cells.Item(2, 18).Value = 4124.067626084168
$ws.Cells.Item(2, 19).Value = 0.01441958321127288
$ws.Cells.Item(2, 20).Value = 0.01441958321127288

$ws.Cells.Item(3, 7).Value = 18.76993966666667
$ws.Cells.Item(3, 8).Value = 56.309819
$ws.Cells.Item(3, 9).Value = 0.1007685501185251
$ws.Cells.Item(3, 10).Value = 0.1007685501185251
$ws.Cells.Item(3, 13).Value = 81.49602766666666
$ws.Cells.Item(3, 15).Value = 0.4776873561738063
$ws.Cells.Item(3, 16).Value = 0.4776873561738064
$ws.Cells.Item(3, 17).Value = 1529.675522376331
$ws.Cells.Item(3, 18).Value = 13767.07970138698
$ws.Cells.Item(3, 19).Value = 0.04813586229158596
$ws.Cells.Item(3, 20).Value = 0.04813586229158597

$ws.Cells.Item(4, 7).Value = 18.76993966666667
$ws.Cells.Item(4, 8).Value = 56.309819
$ws.Cells.Item(4, 9).Value = 0.1007685501185251
$ws.Cells.Item(4, 10).Value = 0.1007685501185251
$ws.Cells.Item(4, 13).Value = 58.29008100000001
$ws.Cells.Item(4, 14).Value = 174.870243
$ws.Cells.Item(4, 15).Value = 0.3416661582321011
$ws.Cells.Item(4, 16).Value = 0.3416661582321012
$ws.Cells.Item(4, 17).Value = 1094.101303535113
$ws.Cells.Item(4, 18).Value = 9846.911731816017
$ws.Cells.Item(4, 19).Value = 0.03442920338961541
$ws.Cells.Item(4, 20).Value = 0.03442920338961542

$ws.Cells.Item(5, 7).Value = 18.76993966666667
$ws.Cells.Item(5, 8).Value = 56.309819
$ws.Cells.Item(5, 9).Value = 0.1007685501185251
$ws.Cells.Item(5, 10).Value = 0.1007685501185251
$ws.Cells.Item(5, 13).Value = 6.406303
$ws.Cells.Item(5, 14).Value = 19.218909
$ws.Cells.Item(5, 15).Value = 0.03755041847481365
$ws.Cells.Item(5, 16).Value = 0.03755041847481365
$ws.Cells.Item(5, 17).Value = 120.2459207963857
$ws.Cells.Item(5, 18).Value = 1082.213287167471
$ws.Cells.Item(5, 19).Value = 0.003783901226050851
$ws.Cells.Item(5, 20).Value = 0.003783901226050851

$ws.Cells.Item(6, 9).Value = 0.5130361557055731
$ws.Cells.Item(6, 10).Value = 0.5130361557055731
$ws.Cells.Item(6, 13).Value = 24.41295733333333
$ws.Cells.Item(6, 14).Value = 73.238872
$ws.Cells.Item(6, 15).Value = 0.1430960671192788
$ws.Cells.Item(6, 16).Value = 0.1430960671192788
$ws.Cells.Item(6, 17).Value = 2332.954300024283
$ws.Cells.Item(6, 18).Value = 20996.58870021854
$ws.Cells.Item(6, 19).Value = 0.07341345617146149
$ws.Cells.Item(6, 20).Value = 0.07341345617146149

$ws.Cells.Item(7, 9).Value = 0.5130361557055731
$ws.Cells.Item(7, 10).Value = 0.5130361557055731
$ws.Cells.Item(7, 13).Value = 81.49602766666666
$ws.Cells.Item(7, 15).Value = 0.4776873561738063
$ws.Cells.Item(7, 16).Value = 0.4776873561738064
$ws.Cells.Item(7, 17).Value = 7787.934316349707
$ws.Cells.Item(7, 18).Value = 70091.40884714737
$ws.Cells.Item(7, 19).Value = 0.2450708848405685
$ws.Cells.Item(7, 20).Value = 0.2450708848405685

$ws.Cells.Item(8, 9).Value = 0.5130361557055731
$ws.Cells.Item(8, 10).Value = 0.5130361557055731
$ws.Cells.Item(8, 13).Value = 58.29008100000001
$ws.Cells.Item(8, 14).Value = 174.870243
$ws.Cells.Item(8, 15).Value = 0.3416661582321011
$ws.Cells.Item(8, 16).Value = 0.3416661582321012
$ws.Cells.Item(8, 17).Value = 5570.324531392855
$ws.Cells.Item(8, 18).Value = 50132.92078253569
$ws.Cells.Item(8, 19).Value = 0.1752870923540892
$ws.Cells.Item(8, 20).Value = 0.1752870923540892

$ws.Cells.Item(9, 9).Value = 0.5130361557055731
$ws.Cells.Item(9, 10).Value = 0.5130361557055731
$ws.Cells.Item(9, 13).Value = 6.406303
$ws.Cells.Item(9, 14).Value = 19.218909
$ws.Cells.Item(9, 15).Value = 0.03755041847481365
$ws.Cells.Item(9, 16).Value = 0.03755041847481365
$ws.Cells.Item(9, 17).Value = 612.1999857306021
$ws.Cells.Item(9, 18).Value = 5509.799871575417
$ws.Cells.Item(9, 19).Value = 0.01926472233945392
$ws.Cells.Item(9, 20).Value = 0.01926472233945392

$ws.Cells.Item(10, 7).Value = 20.061603
$ws.Cells.Item(10, 8).Value = 60.184809
$ws.Cells.Item(10, 9).Value = 0.1077029912330274
$ws.Cells.Item(10, 10).Value = 0.1077029912330274
$ws.Cells.Item(10, 13).Value = 24.41295733333333
$ws.Cells.Item(10, 14).Value = 73.238872
$ws.Cells.Item(10, 15).Value = 0.1430960671192788
$ws.Cells.Item(10, 16).Value = 0.1430960671192788
$ws.Cells.Item(10, 17).Value = 489.7630580772721
$ws.Cells.Item(10, 18).Value = 4407.867522695448
$ws.Cells.Item(10, 19).Value = 0.01541187446242839
$ws.Cells.Item(10, 20).Value = 0.01541187446242839

$ws.Cells.Item(11, 7).Value = 20.061603
$ws.Cells.Item(11, 8).Value = 60.184809
$ws.Cells.Item(11, 9).Value = 0.1077029912330274
$ws.Cells.Item(11, 10).Value = 0.1077029912330274
$ws.Cells.Item(11, 13).Value = 81.49602766666666
$ws.Cells.Item(11, 15).Value = 0.4776873561738063
$ws.Cells.Item(11, 16).Value = 0.4776873561738064
$ws.Cells.Item(11, 17).Value = 1634.940953125683
$ws.Cells.Item(11, 18).Value = 14714.46857813115
$ws.Cells.Item(11, 19).Value = 0.05144835713411552
$ws.Cells.Item(11, 20).Value = 0.05144835713411552

$ws.Cells.Item(12, 7).Value = 20.061603
$ws.Cells.Item(12, 8).Value = 60.184809
$ws.Cells.Item(12, 9).Value = 0.1077029912330274
$ws.Cells.Item(12, 10).Value = 0.1077029912330274
$ws.Cells.Item(12, 13).Value = 58.29008100000001
$ws.Cells.Item(12, 14).Value = 174.870243
$ws.Cells.Item(12, 15).Value = 0.3416661582321011
$ws.Cells.Item(12, 16).Value = 0.3416661582321012
$ws.Cells.Item(12, 17).Value = 1169.392463859843
$ws.Cells.Item(12, 18).Value = 10524.53217473859
$ws.Cells.Item(12, 19).Value = 0.03679846724469416
$ws.Cells.Item(12, 20).Value = 0.03679846724469416

$ws.Cells.Item(13, 7).Value = 20.061603
$ws.Cells.Item(13, 8).Value = 60.184809
$ws.Cells.Item(13, 9).Value = 0.1077029912330274
$ws.Cells.Item(13, 10).Value = 0.1077029912330274
$ws.Cells.Item(13, 13).Value = 6.406303
$ws.Cells.Item(13, 14).Value = 19.218909
$ws.Cells.Item(13, 15).Value = 0.03755041847481365
$ws.Cells.Item(13, 16).Value = 0.03755041847481365
$ws.Cells.Item(13, 17).Value = 128.520707483709
$ws.Cells.Item(13, 18).Value = 1156.686367353381
$ws.Cells.Item(13, 19).Value = 0.004044292391789366
$ws.Cells.Item(13, 20).Value = 0.004044292391789366

$ws.Cells.Item(14, 7).Value = 51.87415833333333
$ws.Cells.Item(14, 8).Value = 155.622475
$ws.Cells.Item(14, 9).Value = 0.2784923029428744
$ws.Cells.Item(14, 10).Value = 0.2784923029428744
$ws.Cells.Item(14, 13).Value = 24.41295733333333
$ws.Cells.Item(14, 14).Value = 73.238872
$ws.Cells.Item(14, 15).Value = 0.1430960671192788
$ws.Cells.Item(14, 16).Value = 0.1430960671192788
$ws.Cells.Item(14, 17).Value = 1266.401614094244
$ws.Cells.Item(14, 18).Value = 11397.6145268482
$ws.Cells.Item(14, 19).Value = 0.03985115327411608
$ws.Cells.Item(14, 20).Value = 0.03985115327411608

$ws.Cells.Item(15, 7).Value = 51.87415833333333
$ws.Cells.Item(15, 8).Value = 155.622475
$ws.Cells.Item(15, 9).Value = 0.2784923029428744
$ws.Cells.Item(15, 10).Value = 0.2784923029428744
$ws.Cells.Item(15, 13).Value = 81.49602766666666
$ws.Cells.Item(15, 15).Value = 0.4776873561738063
$ws.Cells.Item(15, 16).Value = 0.4776873561738064
$ws.Cells.Item(15, 17).Value = 4227.53784271838
$ws.Cells.Item(15, 18).Value = 38047.84058446543
$ws.Cells.Item(15, 19).Value = 0.1330322519075364
$ws.Cells.Item(15, 20).Value = 0.1330322519075364

$ws.Cells.Item(16, 7).Value = 51.87415833333333
$ws.Cells.Item(16, 8).Value = 155.622475
$ws.Cells.Item(16, 9).Value = 0.2784923029428744
$ws.Cells.Item(16, 10).Value = 0.2784923029428744
$ws.Cells.Item(16, 13).Value = 58.29008100000001
$ws.Cells.Item(16, 14).Value = 174.870243
$ws.Cells.Item(16, 15).Value = 0.3416661582321011
$ws.Cells.Item(16, 16).Value = 0.3416661582321012
$ws.Cells.Item(16, 17).Value = 3023.748891056825
$ws.Cells.Item(16, 18).Value = 27213.74001951143
$ws.Cells.Item(16, 19).Value = 0.09515139524370235
$ws.Cells.Item(16, 20).Value = 0.09515139524370236

$ws.Cells.Item(17, 7).Value = 51.87415833333333
$ws.Cells.Item(17, 8).Value = 155.622475
$ws.Cells.Item(17, 9).Value = 0.2784923029428744
$ws.Cells.Item(17, 10).Value = 0.2784923029428744
$ws.Cells.Item(17, 13).Value = 6.406303
$ws.Cells.Item(17, 14).Value = 19.218909
$ws.Cells.Item(17, 15).Value = 0.03755041847481365
$ws.Cells.Item(17, 16).Value = 0.03755041847481365
$ws.Cells.Item(17, 17).Value = 332.3215761533083
$ws.Cells.Item(17, 18).Value = 2990.894185379775
$ws.Cells.Item(17, 19).Value = 0.01045750251751951
$ws.Cells.Item(17, 20).Value = 0.01045750251751951
